$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the unused Sheet2 and Sheet3 tabs
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

$ws = $wb.Worksheets.Item("ICT Data Coordinators")

# Insert a new header row at the top and push existing data down
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Constituency"

# Match the selection recorded in the saved workbook
$ws.Range("C10").Select()
